# Fruta / hortaliza, semanal
# Insert a new weekly data row at row 71 (pushing existing rows 71-90 down to 72-91)
# and populate it with the new observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 71; this shifts rows 71:90 down to 72:91
# and carries down formatting (e.g. the date style in column D) from row 71.
$ws.Rows("71").Insert()

# Fill in the values for the newly inserted row 71.
$ws.Range("A71").Value = 7
$ws.Range("B71").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C71").Value = "Ñuble"
$ws.Range("D71").Value = 44627
$ws.Range("E71").Value = 16
$ws.Range("F71").Value = 100112030
$ws.Range("G71").Value = "Poroto granado"
$ws.Range("H71").Value = "Sin especificar"
$ws.Range("I71").Value = "Primera"
$ws.Range("J71").Value = 60
$ws.Range("K71").Value = 23000
$ws.Range("L71").Value = 24000
$ws.Range("M71").Value = 23500
$ws.Range("N71").Value = "`$/saco 25 kilos"
$ws.Range("O71").Value = "Provincia de Diguillín"
$ws.Range("P71").Value = 940
$ws.Range("Q71").Value = 25
$ws.Range("R71").Value = "Hortaliza"
